# Atualização de bases das ligas, do dia: 02-04-2024 às 23:59
#
# This script corrects mismatched rows of betting-odds data in the
# "Venezuela Primera Division" sheet. Each affected match's data (columns
# B through AC) had been written onto the wrong row; this restores each
# row's data to its correct row (matching row "A" / id stays where it is).
#
# The row-content reassignment (destination row <- source row, copying
# the full B:AC block) is:
#   93  <- 95
#   95  <- 93
#   96  <- 98
#   97  <- 96
#   98  <- 99
#   99  <- 97
#   114 <- 115
#   115 <- 114
#   116 <- 117
#   117 <- 116
#   157 <- 158
#   158 <- 157
#   173 <- 174
#   174 <- 173

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowBlock($row) {
    return $ws.Range("B$row`:AC$row").Value()
}

# Read every affected row's current (pre-edit) B:AC contents first, so
# that later writes never clobber a value before it has been captured.
$rows = @(93, 95, 96, 97, 98, 99, 114, 115, 116, 117, 157, 158, 173, 174)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = Get-RowBlock $r
}

# Destination row -> source row (source's pre-edit data is copied into destination)
$mapping = @{
    93  = 95
    95  = 93
    96  = 98
    97  = 96
    98  = 99
    99  = 97
    114 = 115
    115 = 114
    116 = 117
    117 = 116
    157 = 158
    158 = 157
    173 = 174
    174 = 173
}

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $ws.Range("B$dest`:AC$dest").Value = $snapshot[$src]
}
